$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data adds a new week's worth of readings (2022-01-05, serial 44566)
# for this market/produce combination. In the published sheet this lands at the
# top of the data block (row 48, right after the prior week's last entries),
# pushing all the existing data rows down by two.
$ws.Range("A48:A49").EntireRow.Insert()

# New row 48: "Primera" quality, Región de O'Higgins, $/unidad
$ws.Range("A48").Value = 8
$ws.Range("B48").Value = "Terminal La Palmera de La Serena"
$ws.Range("C48").Value = "Coquimbo"
$ws.Range("D48").Value = 44566
$ws.Range("E48").Value = 4
$ws.Range("F48").Value = 100112028
$ws.Range("G48").Value = "Sandia"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 2000
$ws.Range("K48").Value = 2800
$ws.Range("L48").Value = 3000
$ws.Range("M48").Value = 2900
$ws.Range("N48").Value = "`$/unidad"
$ws.Range("O48").Value = "Región de O'Higgins"
$ws.Range("P48").Value = 2900
$ws.Range("Q48").Value = 1
$ws.Range("R48").Value = "Hortaliza"

# New row 49: "Segunda" quality, Región de O'Higgins, $/unidad
$ws.Range("A49").Value = 8
$ws.Range("B49").Value = "Terminal La Palmera de La Serena"
$ws.Range("C49").Value = "Coquimbo"
$ws.Range("D49").Value = 44566
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 100112028
$ws.Range("G49").Value = "Sandia"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Segunda"
$ws.Range("J49").Value = 1600
$ws.Range("K49").Value = 2000
$ws.Range("L49").Value = 2500
$ws.Range("M49").Value = 2250
$ws.Range("N49").Value = "`$/unidad"
$ws.Range("O49").Value = "Región de O'Higgins"
$ws.Range("P49").Value = 2250
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
